{"js": "// The document body is a single paragraph made of many <w:t> text runs\n// separated by manual line breaks (<w:br/>), all inside one <w:r>.\n// Office.js represents each <w:br/> as a vertical-tab character (\\u000b)\n// inside the paragraph's `text` string. Rebuilding the whole paragraph's\n// text (segments joined by \\u000b, with a trailing \\u000b for the final\n// break) reproduces that exact <w:t>/<w:br/> structure on save.\nconst segments = [\n  \"********************************8\u670830\u65e5\u8bfb\u7ecf\u7ae0\u8282***************************\",\n  \"Chapter 3 of 1_Thessalonians\",\n  \"1.\u6211\u4eec\u65e2\u4e0d\u80fd\u518d\u5fcd\uff0c\u5c31\u613f\u610f\u72ec\u81ea\u7b49\u5728\u96c5\u5178\u3002\",\n  \"2.\u6253\u53d1\u6211\u4eec\u7684\u5144\u5f1f\u5728\u57fa\u7763\u798f\u97f3\u4e0a\u4f5c\u795e\u6267\u4e8b\u7684\u63d0\u6469\u592a\u524d\u53bb\uff0c\uff08\u4f5c\u795e\u6267\u4e8b\u7684\u6709\u53e4\u5377\u4f5c\u4e0e\u795e\u540c\u5de5\u7684\uff09\u575a\u56fa\u4f60\u4eec\uff0c\u5e76\u5728\u4f60\u4eec\u6240\u4fe1\u7684\u9053\u4e0a\u529d\u6170\u4f60\u4eec\u3002\",\n  \"3.\u514d\u5f97\u6709\u4eba\u88ab\u8bf8\u822c\u60a3\u96be\u6447\u52a8\u3002\u56e0\u4e3a\u4f60\u4eec\u81ea\u5df1\u77e5\u9053\u6211\u4eec\u53d7\u60a3\u96be\u539f\u662f\u547d\u5b9a\u7684\u3002\",\n  \"4.\u6211\u4eec\u5728\u4f60\u4eec\u90a3\u91cc\u7684\u65f6\u5019\uff0c\u9884\u5148\u544a\u8bc9\u4f60\u4eec\uff0c\u6211\u4eec\u5fc5\u53d7\u60a3\u96be\uff0c\u4ee5\u540e\u679c\u7136\u5e94\u9a8c\u4e86\uff0c\u4f60\u4eec\u4e5f\u77e5\u9053\u3002\",\n  \"5.\u4e3a\u6b64\uff0c\u6211\u65e2\u4e0d\u80fd\u518d\u5fcd\uff0c\u5c31\u6253\u53d1\u4eba\u53bb\uff0c\u8981\u6653\u5f97\u4f60\u4eec\u7684\u4fe1\u5fc3\u5982\u4f55\uff0c\u6050\u6015\u90a3\u8bf1\u60d1\u4eba\u7684\u5230\u5e95\u8bf1\u60d1\u4e86\u4f60\u4eec\uff0c\u53eb\u6211\u4eec\u7684\u52b3\u82e6\u5f52\u4e8e\u5f92\u7136\u3002\",\n  \"6.\u4f46\u63d0\u6469\u592a\u521a\u624d\u4ece\u4f60\u4eec\u90a3\u91cc\u56de\u6765\uff0c\u5c06\u4f60\u4eec\u4fe1\u5fc3\u548c\u7231\u5fc3\u7684\u597d\u6d88\u606f\u62a5\u7ed9\u6211\u4eec\uff0c\u53c8\u8bf4\u4f60\u4eec\u5e38\u5e38\u8bb0\u5ff5\u6211\u4eec\uff0c\u5207\u5207\u5730\u60f3\u89c1\u6211\u4eec\uff0c\u5982\u540c\u6211\u4eec\u60f3\u89c1\u4f60\u4eec\u4e00\u6837\u3002\",\n  \"7.\u6240\u4ee5\u5f1f\u5144\u4eec\uff0c\u6211\u4eec\u5728\u4e00\u5207\u56f0\u82e6\u60a3\u96be\u4e4b\u4e2d\uff0c\u56e0\u7740\u4f60\u4eec\u7684\u4fe1\u5fc3\u5c31\u5f97\u4e86\u5b89\u6170\u3002\",\n  \"8.\u4f60\u4eec\u82e5\u9760\u4e3b\u7ad9\u7acb\u5f97\u7a33\uff0c\u6211\u4eec\u5c31\u6d3b\u4e86\u3002\",\n  \"9.\u6211\u4eec\u5728\u795e\u9762\u524d\uff0c\u56e0\u7740\u4f60\u4eec\u751a\u662f\u559c\u4e50\uff0c\u4e3a\u8fd9\u4e00\u5207\u559c\u4e50\uff0c\u53ef\u7528\u4f55\u7b49\u7684\u611f\u8c22\uff0c\u4e3a\u4f60\u4eec\u62a5\u7b54\u795e\u5462\uff1f\",\n  \"10.\u6211\u4eec\u663c\u591c\u5207\u5207\u5730\u7948\u6c42\uff0c\u8981\u89c1\u4f60\u4eec\u7684\u9762\uff0c\u8865\u6ee1\u4f60\u4eec\u4fe1\u5fc3\u7684\u4e0d\u8db3\u3002\",\n  \"11.\u613f\u795e\u6211\u4eec\u7684\u7236\uff0c\u548c\u6211\u4eec\u7684\u4e3b\u8036\u7a23\uff0c\u4e00\u76f4\u5f15\u9886\u6211\u4eec\u5230\u4f60\u4eec\u90a3\u91cc\u53bb\u3002\",\n  \"12.\u53c8\u613f\u4e3b\u53eb\u4f60\u4eec\u5f7c\u6b64\u76f8\u7231\u7684\u5fc3\uff0c\u5e76\u7231\u4f17\u4eba\u7684\u5fc3\uff0c\u90fd\u80fd\u589e\u957f\uff0c\u5145\u8db3\uff0c\u5982\u540c\u6211\u4eec\u7231\u4f60\u4eec\u4e00\u6837\u3002\",\n  \"13.\u597d\u4f7f\u4f60\u4eec\uff0c\u5f53\u6211\u4eec\u4e3b\u8036\u7a23\u540c\u4ed6\u4f17\u5723\u5f92\u6765\u7684\u65f6\u5019\uff0c\u5728\u6211\u4eec\u7236\u795e\u9762\u524d\uff0c\u5fc3\u91cc\u575a\u56fa\uff0c\u6210\u4e3a\u5723\u6d01\uff0c\u65e0\u53ef\u8d23\u5907\u3002\",\n  \"Chapter 21 of Proverbs\",\n  \"1.\u738b\u7684\u5fc3\u5728\u8036\u548c\u534e\u624b\u4e2d\uff0c\u597d\u50cf\u9647\u6c9f\u7684\u6c34\uff0c\u968f\u610f\u6d41\u8f6c\u3002\",\n  \"2.\u4eba\u6240\u884c\u7684\uff0c\u5728\u81ea\u5df1\u773c\u4e2d\u90fd\u770b\u4e3a\u6b63\uff0c\u60df\u6709\u8036\u548c\u534e\u8861\u91cf\u4eba\u5fc3\u3002\",\n  \"3.\u884c\u4ec1\u4e49\u516c\u5e73\uff0c\u6bd4\u732e\u796d\u66f4\u8499\u8036\u548c\u534e\u60a6\u7eb3\u3002\",\n  \"4.\u6076\u4eba\u53d1\u8fbe\uff0c\u773c\u9ad8\u5fc3\u50b2\uff0c\u8fd9\u4e43\u662f\u7f6a\u3002\uff08\u53d1\u8fbe\u539f\u6587\u4f5c\u706f\uff09\",\n  \"5.\u6bb7\u52e4\u7b79\u5212\u7684\uff0c\u8db3\u81f4\u4e30\u88d5\u3002\u884c\u4e8b\u6025\u8e81\u7684\uff0c\u90fd\u5fc5\u7f3a\u4e4f\u3002\",\n  \"6.\u7528\u8be1\u8bc8\u4e4b\u820c\u6c42\u8d22\u7684\uff0c\u5c31\u662f\u81ea\u5df1\u53d6\u6b7b\u3002\u6240\u5f97\u4e4b\u8d22\uff0c\u4e43\u662f\u5439\u6765\u5439\u53bb\u7684\u6d6e\u4e91\u3002\",\n  \"7.\u6076\u4eba\u7684\u5f3a\u66b4\uff0c\u5fc5\u5c06\u81ea\u5df1\u626b\u9664\u3002\u56e0\u4ed6\u4eec\u4e0d\u80af\u6309\u516c\u5e73\u884c\u4e8b\u3002\",\n  \"8.\u8d1f\u7f6a\u4e4b\u4eba\u7684\u8def\uff0c\u751a\u662f\u5f2f\u66f2\u3002\u81f3\u4e8e\u6e05\u6d01\u7684\u4eba\uff0c\u4ed6\u6240\u884c\u7684\u4e43\u662f\u6b63\u76f4\u3002\",\n  \"9.\u5b81\u53ef\u4f4f\u5728\u623f\u9876\u7684\u89d2\u4e0a\uff0c\u4e0d\u5728\u5bbd\u9614\u7684\u623f\u5c4b\uff0c\u4e0e\u4e89\u5435\u7684\u5987\u4eba\u540c\u4f4f\u3002\",\n  \"10.\u6076\u4eba\u7684\u5fc3\uff0c\u4e50\u4eba\u53d7\u7978\u3002\u4ed6\u773c\u5e76\u4e0d\u601c\u6064\u90bb\u820d\u3002\",\n  \"11.\u4eb5\u6162\u7684\u4eba\u53d7\u5211\u7f5a\uff0c\u611a\u8499\u7684\u4eba\u5c31\u5f97\u667a\u6167\u3002\u667a\u6167\u4eba\u53d7\u8bad\u8bf2\uff0c\u4fbf\u5f97\u77e5\u8bc6\u3002\",\n  \"12.\u4e49\u4eba\u601d\u60f3\u6076\u4eba\u7684\u5bb6\uff0c\u77e5\u9053\u6076\u4eba\u503e\u5012\uff0c\u5fc5\u81f3\u706d\u4ea1\u3002\",\n  \"13.\u585e\u8033\u4e0d\u542c\u7a77\u4eba\u54c0\u6c42\u7684\uff0c\u4ed6\u5c06\u6765\u547c\u5401\u4e5f\u4e0d\u8499\u5e94\u5141\u3002\",\n  \"14.\u6697\u4e2d\u9001\u7684\u793c\u7269\uff0c\u633d\u56de\u6012\u6c14\u3002\u6000\u4e2d\u640b\u7684\u8d3f\u8d42\uff0c\u6b62\u606f\u66b4\u6012\u3002\",\n  \"15.\u79c9\u516c\u884c\u4e49\uff0c\u4f7f\u4e49\u4eba\u559c\u4e50\uff0c\u4f7f\u4f5c\u5b7d\u7684\u4eba\u8d25\u574f\u3002\",\n  \"16.\u8ff7\u79bb\u901a\u8fbe\u9053\u8def\u7684\uff0c\u5fc5\u4f4f\u5728\u9634\u9b42\u7684\u4f1a\u4e2d\u3002\",\n  \"17.\u7231\u5bb4\u4e50\u7684\uff0c\u5fc5\u81f4\u7a77\u4e4f\u3002\u597d\u9152\u7231\u818f\u6cb9\u7684\uff0c\u5fc5\u4e0d\u5bcc\u8db3\u3002\",\n  \"18.\u6076\u4eba\u4f5c\u4e86\u4e49\u4eba\u7684\u8d4e\u4ef7\u3002\u5978\u8bc8\u4eba\u4ee3\u66ff\u6b63\u76f4\u4eba\u3002\",\n  \"19.\u5b81\u53ef\u4f4f\u5728\u65f7\u91ce\uff0c\u4e0d\u4e0e\u4e89\u5435\u4f7f\u6c14\u7684\u5987\u4eba\u540c\u4f4f\u3002\",\n  \"20.\u667a\u6167\u4eba\u5bb6\u4e2d\u79ef\u84c4\u5b9d\u7269\u818f\u6cb9\u3002\u611a\u6627\u4eba\u968f\u5f97\u6765\u968f\u541e\u4e0b\u3002\",\n  \"21.\u8ffd\u6c42\u516c\u4e49\u4ec1\u6148\u7684\uff0c\u5c31\u5bfb\u5f97\u751f\u547d\uff0c\u516c\u4e49\uff0c\u548c\u5c0a\u8363\u3002\",\n  \"22.\u667a\u6167\u4eba\u722c\u4e0a\u52c7\u58eb\u7684\u57ce\u5899\uff0c\u503e\u8986\u4ed6\u6240\u501a\u9760\u7684\u575a\u5792\u3002\",\n  \"23.\u8c28\u5b88\u53e3\u4e0e\u820c\u7684\uff0c\u5c31\u4fdd\u5b88\u81ea\u5df1\u514d\u53d7\u707e\u96be\u3002\",\n  \"24.\u5fc3\u9a84\u6c14\u50b2\u7684\u4eba\uff0c\u540d\u53eb\u4eb5\u6162\u3002\u4ed6\u884c\u4e8b\u72c2\u5984\uff0c\u90fd\u51fa\u4e8e\u9a84\u50b2\u3002\",\n  \"25.\u61d2\u60f0\u4eba\u7684\u5fc3\u613f\uff0c\u5c06\u4ed6\u6740\u5bb3\uff0c\u56e0\u4e3a\u4ed6\u624b\u4e0d\u80af\u4f5c\u5de5\u3002\",\n  \"26.\u6709\u7ec8\u65e5\u8d2a\u5f97\u65e0\u990d\u7684\uff0c\u4e49\u4eba\u65bd\u820d\u800c\u4e0d\u541d\u60dc\u3002\",\n  \"27.\u6076\u4eba\u7684\u796d\u7269\u662f\u53ef\u618e\u7684\uff0c\u4f55\u51b5\u4ed6\u5b58\u6076\u610f\u6765\u732e\u5462\uff1f\",\n  \"28.\u4f5c\u5047\u89c1\u8bc1\u7684\u5fc5\u706d\u4ea1\uff0c\u60df\u6709\u542c\u771f\u60c5\u800c\u8a00\u7684\uff0c\u5176\u8a00\u957f\u5b58\u3002\",\n  \"29.\u6076\u4eba\u8138\u65e0\u7f9e\u803b\uff0c\u6b63\u76f4\u4eba\u884c\u4e8b\u575a\u5b9a\u3002\",\n  \"30.\u6ca1\u6709\u4eba\u80fd\u4ee5\u667a\u6167\uff0c\u806a\u660e\uff0c\u8c0b\u7565\uff0c\u654c\u6321\u8036\u548c\u534e\u3002\",\n  \"31.\u9a6c\u662f\u4e3a\u6253\u4ed7\u4e4b\u65e5\u9884\u5907\u7684\u3002\u5f97\u80dc\u4e43\u5728\u4e4e\u8036\u548c\u534e\u3002\",\n  \"Chapter 22 of Proverbs\",\n  \"1.\u7f8e\u540d\u80dc\u8fc7\u5927\u8d22\uff0c\u6069\u5ba0\u5f3a\u5982\u91d1\u94f6\u3002\",\n  \"2.\u5bcc\u6237\u7a77\u4eba\uff0c\u5728\u4e16\u76f8\u9047\uff0c\u90fd\u4e3a\u8036\u548c\u534e\u6240\u9020\u3002\",\n  \"3.\u901a\u8fbe\u4eba\u89c1\u7978\u85cf\u8eb2\u3002\u611a\u8499\u4eba\u524d\u5f80\u53d7\u5bb3\u3002\",\n  \"4.\u656c\u754f\u8036\u548c\u534e\u5fc3\u5b58\u8c26\u5351\uff0c\u5c31\u5f97\u5bcc\u6709\uff0c\u5c0a\u8363\uff0c\u751f\u547d\uff0c\u4e3a\u8d4f\u8d50\u3002\",\n  \"5.\u4e56\u50fb\u4eba\u7684\u8def\u4e0a\uff0c\u6709\u8346\u68d8\u548c\u7f51\u7f57\u3002\u4fdd\u5b88\u81ea\u5df1\u751f\u547d\u7684\u3002\u5fc5\u8981\u8fdc\u79bb\u3002\",\n  \"6.\u6559\u517b\u5b69\u7ae5\uff0c\u4f7f\u4ed6\u8d70\u5f53\u884c\u7684\u9053\uff0c\u5c31\u662f\u5230\u8001\u4ed6\u4e5f\u4e0d\u504f\u79bb\u3002\",\n  \"7.\u5bcc\u6237\u7ba1\u8f96\u7a77\u4eba\uff0c\u6b20\u503a\u7684\u662f\u503a\u4e3b\u7684\u4ec6\u4eba\u3002\",\n  \"8.\u6492\u7f6a\u5b7d\u7684\uff0c\u5fc5\u6536\u707e\u7978\u3002\u4ed6\u901e\u6012\u7684\u6756\uff0c\u4e5f\u5fc5\u5e9f\u6389\u3002\",\n  \"9.\u773c\u76ee\u6148\u5584\u7684\uff0c\u5c31\u5fc5\u8499\u798f\u3002\u56e0\u4ed6\u5c06\u98df\u7269\u5206\u7ed9\u7a77\u4eba\u3002\",\n  \"10.\u8d76\u51fa\u4eb5\u6162\u4eba\uff0c\u4e89\u7aef\u5c31\u6d88\u9664\uff0c\u5206\u4e89\u548c\u7f9e\u8fb1\uff0c\u4e5f\u5fc5\u6b62\u606f\u3002\",\n  \"11.\u559c\u7231\u6e05\u5fc3\u7684\u4eba\uff0c\u56e0\u4ed6\u5634\u4e0a\u7684\u6069\u8a00\uff0c\u738b\u5fc5\u4e0e\u4ed6\u4e3a\u53cb\u3002\",\n  \"12.\u8036\u548c\u534e\u7684\u773c\u76ee\uff0c\u7737\u987e\u806a\u660e\u4eba\u3002\u5374\u503e\u8d25\u5978\u8bc8\u4eba\u7684\u8a00\u8bed\u3002\",\n  \"13.\u61d2\u60f0\u4eba\u8bf4\uff0c\u5916\u5934\u6709\u72ee\u5b50\uff0c\u6211\u5728\u8857\u4e0a\uff0c\u5c31\u5fc5\u88ab\u6740\u3002\",\n  \"14.\u6deb\u5987\u7684\u53e3\u4e3a\u6df1\u5751\uff0c\u8036\u548c\u534e\u6240\u618e\u6076\u7684\uff0c\u5fc5\u9677\u5728\u5176\u4e2d\u3002\",\n  \"15.\u611a\u8499\u8ff7\u4f4f\u5b69\u7ae5\u7684\u5fc3\uff0c\u7528\u7ba1\u6559\u7684\u6756\u53ef\u4ee5\u8fdc\u8fdc\u8d76\u9664\u3002\",\n  \"16.\u6b3a\u538b\u8d2b\u7a77\u4e3a\u8981\u5229\u5df1\u7684\uff0c\u5e76\u9001\u793c\u4e0e\u5bcc\u6237\u7684\uff0c\u90fd\u5fc5\u7f3a\u4e4f\u3002\",\n  \"17.\u4f60\u987b\u4fa7\u8033\u542c\u53d7\u667a\u6167\u4eba\u7684\u8a00\u8bed\uff0c\u7559\u5fc3\u9886\u4f1a\u6211\u7684\u77e5\u8bc6\uff0c\",\n  \"18.\u4f60\u82e5\u5fc3\u4e2d\u5b58\u8bb0\uff0c\u5634\u4e0a\u54ac\u5b9a\uff0c\u8fd9\u4fbf\u4e3a\u7f8e\u3002\",\n  \"19.\u6211\u4eca\u65e5\u4ee5\u6b64\u7279\u7279\u6307\u6559\u4f60\uff0c\u4e3a\u8981\u4f7f\u4f60\u501a\u9760\u8036\u548c\u534e\u3002\",\n  \"20.\u8c0b\u7565\u548c\u77e5\u8bc6\u7684\u7f8e\u4e8b\uff0c\u6211\u5c82\u6ca1\u6709\u5199\u7ed9\u4f60\u5417\uff1f\",\n  \"21.\u8981\u4f7f\u4f60\u77e5\u9053\u771f\u8a00\u7684\u5b9e\u7406\uff0c\u4f60\u597d\u5c06\u771f\u8a00\u56de\u8986\u90a3\u6253\u53d1\u4f60\u6765\u7684\u4eba\u3002\",\n  \"22.\u8d2b\u7a77\u4eba\uff0c\u4f60\u4e0d\u53ef\u56e0\u4ed6\u8d2b\u7a77\uff0c\u5c31\u62a2\u593a\u4ed6\u7684\u7269\u3002\u4e5f\u4e0d\u53ef\u5728\u57ce\u95e8\u53e3\u6b3a\u538b\u56f0\u82e6\u4eba\u3002\",\n  \"23.\u56e0\u8036\u548c\u534e\u5fc5\u4e3a\u4ed6\u8fa8\u5c48\u3002\u62a2\u593a\u4ed6\u7684\uff0c\u8036\u548c\u534e\u5fc5\u593a\u53d6\u90a3\u4eba\u7684\u547d\u3002\",\n  \"24.\u597d\u751f\u6c14\u7684\u4eba\uff0c\u4e0d\u53ef\u4e0e\u4ed6\u7ed3\u4ea4\u3002\u66b4\u6012\u7684\u4eba\uff0c\u4e0d\u53ef\u4e0e\u4ed6\u6765\u5f80\u3002\",\n  \"25.\u6050\u6015\u4f60\u6548\u6cd5\u4ed6\u7684\u884c\u4e3a\uff0c\u81ea\u5df1\u5c31\u9677\u5728\u7f51\u7f57\u91cc\u3002\",\n  \"26.\u4e0d\u8981\u4e0e\u4eba\u51fb\u638c\uff0c\u4e0d\u8981\u4e3a\u6b20\u503a\u7684\u4f5c\u4fdd\u3002\",\n  \"27.\u4f60\u82e5\u6ca1\u6709\u4ec0\u4e48\u507f\u8fd8\uff0c\u4f55\u5fc5\u4f7f\u4eba\u593a\u53bb\u4f60\u7761\u5367\u7684\u5e8a\u5462\uff1f\",\n  \"28.\u4f60\u5148\u7956\u6240\u7acb\u7684\u5730\u754c\uff0c\u4f60\u4e0d\u53ef\u632a\u79fb\u3002\",\n  \"29.\u4f60\u770b\u89c1\u529e\u4e8b\u6bb7\u52e4\u7684\u4eba\u5417\uff1f\u4ed6\u5fc5\u7ad9\u5728\u541b\u738b\u9762\u524d\uff0c\u5fc5\u4e0d\u7ad9\u5728\u4e0b\u8d31\u4eba\u9762\u524d\u3002\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst newText = segments.join(\"\\u000b\") + \"\\u000b\";\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# The document body is one paragraph made of many <w:t> text segments\n# separated by manual line breaks (<w:br/>), all inside a single <w:r>.\n# In the Word object model, Range.Text exposes each <w:br/> as a vertical\n# tab character (Chr(11)). Rebuilding the whole paragraph's Range.Text from\n# the segments joined by Chr(11), with a trailing Chr(11) for the final\n# break, reproduces that exact <w:t>/<w:br/> run structure on save.\n\n$segments = @(\n    '********************************8\u670830\u65e5\u8bfb\u7ecf\u7ae0\u8282***************************',\n    'Chapter 3 of 1_Thessalonians',\n    '1.\u6211\u4eec\u65e2\u4e0d\u80fd\u518d\u5fcd\uff0c\u5c31\u613f\u610f\u72ec\u81ea\u7b49\u5728\u96c5\u5178\u3002',\n    '2.\u6253\u53d1\u6211\u4eec\u7684\u5144\u5f1f\u5728\u57fa\u7763\u798f\u97f3\u4e0a\u4f5c\u795e\u6267\u4e8b\u7684\u63d0\u6469\u592a\u524d\u53bb\uff0c\uff08\u4f5c\u795e\u6267\u4e8b\u7684\u6709\u53e4\u5377\u4f5c\u4e0e\u795e\u540c\u5de5\u7684\uff09\u575a\u56fa\u4f60\u4eec\uff0c\u5e76\u5728\u4f60\u4eec\u6240\u4fe1\u7684\u9053\u4e0a\u529d\u6170\u4f60\u4eec\u3002',\n    '3.\u514d\u5f97\u6709\u4eba\u88ab\u8bf8\u822c\u60a3\u96be\u6447\u52a8\u3002\u56e0\u4e3a\u4f60\u4eec\u81ea\u5df1\u77e5\u9053\u6211\u4eec\u53d7\u60a3\u96be\u539f\u662f\u547d\u5b9a\u7684\u3002',\n    '4.\u6211\u4eec\u5728\u4f60\u4eec\u90a3\u91cc\u7684\u65f6\u5019\uff0c\u9884\u5148\u544a\u8bc9\u4f60\u4eec\uff0c\u6211\u4eec\u5fc5\u53d7\u60a3\u96be\uff0c\u4ee5\u540e\u679c\u7136\u5e94\u9a8c\u4e86\uff0c\u4f60\u4eec\u4e5f\u77e5\u9053\u3002',\n    '5.\u4e3a\u6b64\uff0c\u6211\u65e2\u4e0d\u80fd\u518d\u5fcd\uff0c\u5c31\u6253\u53d1\u4eba\u53bb\uff0c\u8981\u6653\u5f97\u4f60\u4eec\u7684\u4fe1\u5fc3\u5982\u4f55\uff0c\u6050\u6015\u90a3\u8bf1\u60d1\u4eba\u7684\u5230\u5e95\u8bf1\u60d1\u4e86\u4f60\u4eec\uff0c\u53eb\u6211\u4eec\u7684\u52b3\u82e6\u5f52\u4e8e\u5f92\u7136\u3002',\n    '6.\u4f46\u63d0\u6469\u592a\u521a\u624d\u4ece\u4f60\u4eec\u90a3\u91cc\u56de\u6765\uff0c\u5c06\u4f60\u4eec\u4fe1\u5fc3\u548c\u7231\u5fc3\u7684\u597d\u6d88\u606f\u62a5\u7ed9\u6211\u4eec\uff0c\u53c8\u8bf4\u4f60\u4eec\u5e38\u5e38\u8bb0\u5ff5\u6211\u4eec\uff0c\u5207\u5207\u5730\u60f3\u89c1\u6211\u4eec\uff0c\u5982\u540c\u6211\u4eec\u60f3\u89c1\u4f60\u4eec\u4e00\u6837\u3002',\n    '7.\u6240\u4ee5\u5f1f\u5144\u4eec\uff0c\u6211\u4eec\u5728\u4e00\u5207\u56f0\u82e6\u60a3\u96be\u4e4b\u4e2d\uff0c\u56e0\u7740\u4f60\u4eec\u7684\u4fe1\u5fc3\u5c31\u5f97\u4e86\u5b89\u6170\u3002',\n    '8.\u4f60\u4eec\u82e5\u9760\u4e3b\u7ad9\u7acb\u5f97\u7a33\uff0c\u6211\u4eec\u5c31\u6d3b\u4e86\u3002',\n    '9.\u6211\u4eec\u5728\u795e\u9762\u524d\uff0c\u56e0\u7740\u4f60\u4eec\u751a\u662f\u559c\u4e50\uff0c\u4e3a\u8fd9\u4e00\u5207\u559c\u4e50\uff0c\u53ef\u7528\u4f55\u7b49\u7684\u611f\u8c22\uff0c\u4e3a\u4f60\u4eec\u62a5\u7b54\u795e\u5462\uff1f',\n    '10.\u6211\u4eec\u663c\u591c\u5207\u5207\u5730\u7948\u6c42\uff0c\u8981\u89c1\u4f60\u4eec\u7684\u9762\uff0c\u8865\u6ee1\u4f60\u4eec\u4fe1\u5fc3\u7684\u4e0d\u8db3\u3002',\n    '11.\u613f\u795e\u6211\u4eec\u7684\u7236\uff0c\u548c\u6211\u4eec\u7684\u4e3b\u8036\u7a23\uff0c\u4e00\u76f4\u5f15\u9886\u6211\u4eec\u5230\u4f60\u4eec\u90a3\u91cc\u53bb\u3002',\n    '12.\u53c8\u613f\u4e3b\u53eb\u4f60\u4eec\u5f7c\u6b64\u76f8\u7231\u7684\u5fc3\uff0c\u5e76\u7231\u4f17\u4eba\u7684\u5fc3\uff0c\u90fd\u80fd\u589e\u957f\uff0c\u5145\u8db3\uff0c\u5982\u540c\u6211\u4eec\u7231\u4f60\u4eec\u4e00\u6837\u3002',\n    '13.\u597d\u4f7f\u4f60\u4eec\uff0c\u5f53\u6211\u4eec\u4e3b\u8036\u7a23\u540c\u4ed6\u4f17\u5723\u5f92\u6765\u7684\u65f6\u5019\uff0c\u5728\u6211\u4eec\u7236\u795e\u9762\u524d\uff0c\u5fc3\u91cc\u575a\u56fa\uff0c\u6210\u4e3a\u5723\u6d01\uff0c\u65e0\u53ef\u8d23\u5907\u3002',\n    'Chapter 21 of Proverbs',\n    '1.\u738b\u7684\u5fc3\u5728\u8036\u548c\u534e\u624b\u4e2d\uff0c\u597d\u50cf\u9647\u6c9f\u7684\u6c34\uff0c\u968f\u610f\u6d41\u8f6c\u3002',\n    '2.\u4eba\u6240\u884c\u7684\uff0c\u5728\u81ea\u5df1\u773c\u4e2d\u90fd\u770b\u4e3a\u6b63\uff0c\u60df\u6709\u8036\u548c\u534e\u8861\u91cf\u4eba\u5fc3\u3002',\n    '3.\u884c\u4ec1\u4e49\u516c\u5e73\uff0c\u6bd4\u732e\u796d\u66f4\u8499\u8036\u548c\u534e\u60a6\u7eb3\u3002',\n    '4.\u6076\u4eba\u53d1\u8fbe\uff0c\u773c\u9ad8\u5fc3\u50b2\uff0c\u8fd9\u4e43\u662f\u7f6a\u3002\uff08\u53d1\u8fbe\u539f\u6587\u4f5c\u706f\uff09',\n    '5.\u6bb7\u52e4\u7b79\u5212\u7684\uff0c\u8db3\u81f4\u4e30\u88d5\u3002\u884c\u4e8b\u6025\u8e81\u7684\uff0c\u90fd\u5fc5\u7f3a\u4e4f\u3002',\n    '6.\u7528\u8be1\u8bc8\u4e4b\u820c\u6c42\u8d22\u7684\uff0c\u5c31\u662f\u81ea\u5df1\u53d6\u6b7b\u3002\u6240\u5f97\u4e4b\u8d22\uff0c\u4e43\u662f\u5439\u6765\u5439\u53bb\u7684\u6d6e\u4e91\u3002',\n    '7.\u6076\u4eba\u7684\u5f3a\u66b4\uff0c\u5fc5\u5c06\u81ea\u5df1\u626b\u9664\u3002\u56e0\u4ed6\u4eec\u4e0d\u80af\u6309\u516c\u5e73\u884c\u4e8b\u3002',\n    '8.\u8d1f\u7f6a\u4e4b\u4eba\u7684\u8def\uff0c\u751a\u662f\u5f2f\u66f2\u3002\u81f3\u4e8e\u6e05\u6d01\u7684\u4eba\uff0c\u4ed6\u6240\u884c\u7684\u4e43\u662f\u6b63\u76f4\u3002',\n    '9.\u5b81\u53ef\u4f4f\u5728\u623f\u9876\u7684\u89d2\u4e0a\uff0c\u4e0d\u5728\u5bbd\u9614\u7684\u623f\u5c4b\uff0c\u4e0e\u4e89\u5435\u7684\u5987\u4eba\u540c\u4f4f\u3002',\n    '10.\u6076\u4eba\u7684\u5fc3\uff0c\u4e50\u4eba\u53d7\u7978\u3002\u4ed6\u773c\u5e76\u4e0d\u601c\u6064\u90bb\u820d\u3002',\n    '11.\u4eb5\u6162\u7684\u4eba\u53d7\u5211\u7f5a\uff0c\u611a\u8499\u7684\u4eba\u5c31\u5f97\u667a\u6167\u3002\u667a\u6167\u4eba\u53d7\u8bad\u8bf2\uff0c\u4fbf\u5f97\u77e5\u8bc6\u3002',\n    '12.\u4e49\u4eba\u601d\u60f3\u6076\u4eba\u7684\u5bb6\uff0c\u77e5\u9053\u6076\u4eba\u503e\u5012\uff0c\u5fc5\u81f3\u706d\u4ea1\u3002',\n    '13.\u585e\u8033\u4e0d\u542c\u7a77\u4eba\u54c0\u6c42\u7684\uff0c\u4ed6\u5c06\u6765\u547c\u5401\u4e5f\u4e0d\u8499\u5e94\u5141\u3002',\n    '14.\u6697\u4e2d\u9001\u7684\u793c\u7269\uff0c\u633d\u56de\u6012\u6c14\u3002\u6000\u4e2d\u640b\u7684\u8d3f\u8d42\uff0c\u6b62\u606f\u66b4\u6012\u3002',\n    '15.\u79c9\u516c\u884c\u4e49\uff0c\u4f7f\u4e49\u4eba\u559c\u4e50\uff0c\u4f7f\u4f5c\u5b7d\u7684\u4eba\u8d25\u574f\u3002',\n    '16.\u8ff7\u79bb\u901a\u8fbe\u9053\u8def\u7684\uff0c\u5fc5\u4f4f\u5728\u9634\u9b42\u7684\u4f1a\u4e2d\u3002',\n    '17.\u7231\u5bb4\u4e50\u7684\uff0c\u5fc5\u81f4\u7a77\u4e4f\u3002\u597d\u9152\u7231\u818f\u6cb9\u7684\uff0c\u5fc5\u4e0d\u5bcc\u8db3\u3002',\n    '18.\u6076\u4eba\u4f5c\u4e86\u4e49\u4eba\u7684\u8d4e\u4ef7\u3002\u5978\u8bc8\u4eba\u4ee3\u66ff\u6b63\u76f4\u4eba\u3002',\n    '19.\u5b81\u53ef\u4f4f\u5728\u65f7\u91ce\uff0c\u4e0d\u4e0e\u4e89\u5435\u4f7f\u6c14\u7684\u5987\u4eba\u540c\u4f4f\u3002',\n    '20.\u667a\u6167\u4eba\u5bb6\u4e2d\u79ef\u84c4\u5b9d\u7269\u818f\u6cb9\u3002\u611a\u6627\u4eba\u968f\u5f97\u6765\u968f\u541e\u4e0b\u3002',\n    '21.\u8ffd\u6c42\u516c\u4e49\u4ec1\u6148\u7684\uff0c\u5c31\u5bfb\u5f97\u751f\u547d\uff0c\u516c\u4e49\uff0c\u548c\u5c0a\u8363\u3002',\n    '22.\u667a\u6167\u4eba\u722c\u4e0a\u52c7\u58eb\u7684\u57ce\u5899\uff0c\u503e\u8986\u4ed6\u6240\u501a\u9760\u7684\u575a\u5792\u3002',\n    '23.\u8c28\u5b88\u53e3\u4e0e\u820c\u7684\uff0c\u5c31\u4fdd\u5b88\u81ea\u5df1\u514d\u53d7\u707e\u96be\u3002',\n    '24.\u5fc3\u9a84\u6c14\u50b2\u7684\u4eba\uff0c\u540d\u53eb\u4eb5\u6162\u3002\u4ed6\u884c\u4e8b\u72c2\u5984\uff0c\u90fd\u51fa\u4e8e\u9a84\u50b2\u3002',\n    '25.\u61d2\u60f0\u4eba\u7684\u5fc3\u613f\uff0c\u5c06\u4ed6\u6740\u5bb3\uff0c\u56e0\u4e3a\u4ed6\u624b\u4e0d\u80af\u4f5c\u5de5\u3002',\n    '26.\u6709\u7ec8\u65e5\u8d2a\u5f97\u65e0\u990d\u7684\uff0c\u4e49\u4eba\u65bd\u820d\u800c\u4e0d\u541d\u60dc\u3002',\n    '27.\u6076\u4eba\u7684\u796d\u7269\u662f\u53ef\u618e\u7684\uff0c\u4f55\u51b5\u4ed6\u5b58\u6076\u610f\u6765\u732e\u5462\uff1f',\n    '28.\u4f5c\u5047\u89c1\u8bc1\u7684\u5fc5\u706d\u4ea1\uff0c\u60df\u6709\u542c\u771f\u60c5\u800c\u8a00\u7684\uff0c\u5176\u8a00\u957f\u5b58\u3002',\n    '29.\u6076\u4eba\u8138\u65e0\u7f9e\u803b\uff0c\u6b63\u76f4\u4eba\u884c\u4e8b\u575a\u5b9a\u3002',\n    '30.\u6ca1\u6709\u4eba\u80fd\u4ee5\u667a\u6167\uff0c\u806a\u660e\uff0c\u8c0b\u7565\uff0c\u654c\u6321\u8036\u548c\u534e\u3002',\n    '31.\u9a6c\u662f\u4e3a\u6253\u4ed7\u4e4b\u65e5\u9884\u5907\u7684\u3002\u5f97\u80dc\u4e43\u5728\u4e4e\u8036\u548c\u534e\u3002',\n    'Chapter 22 of Proverbs',\n    '1.\u7f8e\u540d\u80dc\u8fc7\u5927\u8d22\uff0c\u6069\u5ba0\u5f3a\u5982\u91d1\u94f6\u3002',\n    '2.\u5bcc\u6237\u7a77\u4eba\uff0c\u5728\u4e16\u76f8\u9047\uff0c\u90fd\u4e3a\u8036\u548c\u534e\u6240\u9020\u3002',\n    '3.\u901a\u8fbe\u4eba\u89c1\u7978\u85cf\u8eb2\u3002\u611a\u8499\u4eba\u524d\u5f80\u53d7\u5bb3\u3002',\n    '4.\u656c\u754f\u8036\u548c\u534e\u5fc3\u5b58\u8c26\u5351\uff0c\u5c31\u5f97\u5bcc\u6709\uff0c\u5c0a\u8363\uff0c\u751f\u547d\uff0c\u4e3a\u8d4f\u8d50\u3002',\n    '5.\u4e56\u50fb\u4eba\u7684\u8def\u4e0a\uff0c\u6709\u8346\u68d8\u548c\u7f51\u7f57\u3002\u4fdd\u5b88\u81ea\u5df1\u751f\u547d\u7684\u3002\u5fc5\u8981\u8fdc\u79bb\u3002',\n    '6.\u6559\u517b\u5b69\u7ae5\uff0c\u4f7f\u4ed6\u8d70\u5f53\u884c\u7684\u9053\uff0c\u5c31\u662f\u5230\u8001\u4ed6\u4e5f\u4e0d\u504f\u79bb\u3002',\n    '7.\u5bcc\u6237\u7ba1\u8f96\u7a77\u4eba\uff0c\u6b20\u503a\u7684\u662f\u503a\u4e3b\u7684\u4ec6\u4eba\u3002',\n    '8.\u6492\u7f6a\u5b7d\u7684\uff0c\u5fc5\u6536\u707e\u7978\u3002\u4ed6\u901e\u6012\u7684\u6756\uff0c\u4e5f\u5fc5\u5e9f\u6389\u3002',\n    '9.\u773c\u76ee\u6148\u5584\u7684\uff0c\u5c31\u5fc5\u8499\u798f\u3002\u56e0\u4ed6\u5c06\u98df\u7269\u5206\u7ed9\u7a77\u4eba\u3002',\n    '10.\u8d76\u51fa\u4eb5\u6162\u4eba\uff0c\u4e89\u7aef\u5c31\u6d88\u9664\uff0c\u5206\u4e89\u548c\u7f9e\u8fb1\uff0c\u4e5f\u5fc5\u6b62\u606f\u3002',\n    '11.\u559c\u7231\u6e05\u5fc3\u7684\u4eba\uff0c\u56e0\u4ed6\u5634\u4e0a\u7684\u6069\u8a00\uff0c\u738b\u5fc5\u4e0e\u4ed6\u4e3a\u53cb\u3002',\n    '12.\u8036\u548c\u534e\u7684\u773c\u76ee\uff0c\u7737\u987e\u806a\u660e\u4eba\u3002\u5374\u503e\u8d25\u5978\u8bc8\u4eba\u7684\u8a00\u8bed\u3002',\n    '13.\u61d2\u60f0\u4eba\u8bf4\uff0c\u5916\u5934\u6709\u72ee\u5b50\uff0c\u6211\u5728\u8857\u4e0a\uff0c\u5c31\u5fc5\u88ab\u6740\u3002',\n    '14.\u6deb\u5987\u7684\u53e3\u4e3a\u6df1\u5751\uff0c\u8036\u548c\u534e\u6240\u618e\u6076\u7684\uff0c\u5fc5\u9677\u5728\u5176\u4e2d\u3002',\n    '15.\u611a\u8499\u8ff7\u4f4f\u5b69\u7ae5\u7684\u5fc3\uff0c\u7528\u7ba1\u6559\u7684\u6756\u53ef\u4ee5\u8fdc\u8fdc\u8d76\u9664\u3002',\n    '16.\u6b3a\u538b\u8d2b\u7a77\u4e3a\u8981\u5229\u5df1\u7684\uff0c\u5e76\u9001\u793c\u4e0e\u5bcc\u6237\u7684\uff0c\u90fd\u5fc5\u7f3a\u4e4f\u3002',\n    '17.\u4f60\u987b\u4fa7\u8033\u542c\u53d7\u667a\u6167\u4eba\u7684\u8a00\u8bed\uff0c\u7559\u5fc3\u9886\u4f1a\u6211\u7684\u77e5\u8bc6\uff0c',\n    '18.\u4f60\u82e5\u5fc3\u4e2d\u5b58\u8bb0\uff0c\u5634\u4e0a\u54ac\u5b9a\uff0c\u8fd9\u4fbf\u4e3a\u7f8e\u3002',\n    '19.\u6211\u4eca\u65e5\u4ee5\u6b64\u7279\u7279\u6307\u6559\u4f60\uff0c\u4e3a\u8981\u4f7f\u4f60\u501a\u9760\u8036\u548c\u534e\u3002',\n    '20.\u8c0b\u7565\u548c\u77e5\u8bc6\u7684\u7f8e\u4e8b\uff0c\u6211\u5c82\u6ca1\u6709\u5199\u7ed9\u4f60\u5417\uff1f',\n    '21.\u8981\u4f7f\u4f60\u77e5\u9053\u771f\u8a00\u7684\u5b9e\u7406\uff0c\u4f60\u597d\u5c06\u771f\u8a00\u56de\u8986\u90a3\u6253\u53d1\u4f60\u6765\u7684\u4eba\u3002',\n    '22.\u8d2b\u7a77\u4eba\uff0c\u4f60\u4e0d\u53ef\u56e0\u4ed6\u8d2b\u7a77\uff0c\u5c31\u62a2\u593a\u4ed6\u7684\u7269\u3002\u4e5f\u4e0d\u53ef\u5728\u57ce\u95e8\u53e3\u6b3a\u538b\u56f0\u82e6\u4eba\u3002',\n    '23.\u56e0\u8036\u548c\u534e\u5fc5\u4e3a\u4ed6\u8fa8\u5c48\u3002\u62a2\u593a\u4ed6\u7684\uff0c\u8036\u548c\u534e\u5fc5\u593a\u53d6\u90a3\u4eba\u7684\u547d\u3002',\n    '24.\u597d\u751f\u6c14\u7684\u4eba\uff0c\u4e0d\u53ef\u4e0e\u4ed6\u7ed3\u4ea4\u3002\u66b4\u6012\u7684\u4eba\uff0c\u4e0d\u53ef\u4e0e\u4ed6\u6765\u5f80\u3002',\n    '25.\u6050\u6015\u4f60\u6548\u6cd5\u4ed6\u7684\u884c\u4e3a\uff0c\u81ea\u5df1\u5c31\u9677\u5728\u7f51\u7f57\u91cc\u3002',\n    '26.\u4e0d\u8981\u4e0e\u4eba\u51fb\u638c\uff0c\u4e0d\u8981\u4e3a\u6b20\u503a\u7684\u4f5c\u4fdd\u3002',\n    '27.\u4f60\u82e5\u6ca1\u6709\u4ec0\u4e48\u507f\u8fd8\uff0c\u4f55\u5fc5\u4f7f\u4eba\u593a\u53bb\u4f60\u7761\u5367\u7684\u5e8a\u5462\uff1f',\n    '28.\u4f60\u5148\u7956\u6240\u7acb\u7684\u5730\u754c\uff0c\u4f60\u4e0d\u53ef\u632a\u79fb\u3002',\n    '29.\u4f60\u770b\u89c1\u529e\u4e8b\u6bb7\u52e4\u7684\u4eba\u5417\uff1f\u4ed6\u5fc5\u7ad9\u5728\u541b\u738b\u9762\u524d\uff0c\u5fc5\u4e0d\u7ad9\u5728\u4e0b\u8d31\u4eba\u9762\u524d\u3002'\n)\n\n$d = $word.ActiveDocument\n\n$br = [char]11\n$newText = ($segments -join $br) + $br\n\n$firstParagraph = $d.Paragraphs.First\n$firstParagraph.Range.Text = $newText\n"}
